# chore: remove invoice feature; fix monthly export weekly sheets + admin cache-bust
#
# Updates the "Weekly Timesheet" sheet for Doug Kinsey's 2026-01-26 week:
# client names in column B and the Hours/Rate/Total figures (and the
# SUBTOTAL row) in columns C/E/F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Client names (column B) ---
$ws.Range("B2").Value = "McGill"
$ws.Range("B3").Value = "Funke (Maint Items)"
$ws.Range("B4").Value = "Caputo (Maint Items)"
$ws.Range("B5").Value = "TOTAL:"
$ws.Range("B6").Value = "Lynn"
$ws.Range("B7").Value = "Lynn"

# --- Hours (column C) ---
$ws.Range("C2").Value = 10.5
$ws.Range("C3").Value = 10.5
$ws.Range("C4").Value = 8.5
$ws.Range("C5").Value = 7.5
$ws.Range("C6").Value = 3
$ws.Range("C7").Value = 5.5

# --- Rate (column E) ---
$ws.Range("E2").Value = 65
$ws.Range("E3").Value = 65
$ws.Range("E4").Value = 65
$ws.Range("E5").Value = 65
$ws.Range("E6").Value = 65
$ws.Range("E7").Value = 65

# --- Total (column F) ---
$ws.Range("F2").Value = 682.5
$ws.Range("F3").Value = 682.5
$ws.Range("F4").Value = 552.5
$ws.Range("F5").Value = 487.5
$ws.Range("F6").Value = 195
$ws.Range("F7").Value = 536.25

# --- SUBTOTAL row (row 9) ---
$ws.Range("C9").Value = 45.5
$ws.Range("D9").Value = "Reg: 40 / OT: 5.5"
$ws.Range("F9").Value = 3136.25
